$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data (post-edit) -------------------------------------------------
# Column B: updated quantities (rows 1-14)
$B = @(699,312,428,495,464,584,848,760,897,1163,1670,1770,786,2224)
# New column C: percentages (rows 1-14)
$C = @(0.05,0.02,0.03,0.03,0.03,0.04,0.06,0.05,0.06,0.08,0.12,0.13,0.06,0.16)
# Old column C values (unchanged), now living in column D (rows 1-13 only)
$D = @(0,5,15,25,35,45,55,65,75,85,95,100,100)

# --- Shift the existing C/D/E columns one slot to the right (rows 1-13) ---
# so a brand-new column C can be inserted for the percentages, while leaving
# row 15 untouched (it never had data past column B before this edit).

# old E (style-only, empty, s=1) -> F
$ws.Range("E1:E13").Copy($ws.Range("F1:F13")) | Out-Null

# old D (formula "=C*B") -> E, formula rewritten to reference the new column D
for ($r = 1; $r -le 13; $r++) {
    $ws.Cells.Item($r, 5).ClearFormats() | Out-Null
    $ws.Cells.Item($r, 5).Formula = "=D$r*B$r"
}

# old C (plain numeric values) -> D
for ($r = 1; $r -le 13; $r++) {
    $ws.Cells.Item($r, 4).Value = $D[$r - 1]
}

# --- Row 14: move the average formula (old E14, style idx 2) to F14 -------
$ws.Range("E14").Copy($ws.Range("F14")) | Out-Null
$ws.Range("F14").Formula = "=SUM(E1:E13)/SUM(B1:B13)"
$ws.Range("E14").Clear() | Out-Null

# --- Row 15: new running total for the percentage column ------------------
# (written before the column-C percent formatting loop below so the newly
# formatted column doesn't bleed its number format into this cell)
$ws.Range("C15").Formula = "=SUM(C1:C14)"

# --- Column B: write the updated quantities --------------------------------
for ($r = 1; $r -le 14; $r++) {
    $ws.Cells.Item($r, 2).Value = $B[$r - 1]
}

# --- New column C: percentages with a 0% number format ---------------------
for ($r = 1; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).NumberFormat = "0%"
    $ws.Cells.Item($r, 3).Value = $C[$r - 1]
}

# --- Cosmetic: selection + page setup (paper size / orientation) ----------
$ws.Range("C1:C14").Select() | Out-Null

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
